# Registration Test script added
# Inserts a new "Registration Details" block (rows 2-6) above the existing
# login-object table on Sheet1, pushing the old data down to rows 7-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room: insert 5 blank rows before the current row 2 -------------
for ($i = 0; $i -lt 5; $i++) {
    $ws.Rows("2:2").Insert()
}

# --- New "Registration Details" section -----------------------------------
$ws.Range("A3").Value = "Registration Details"

$ws.Range("A4").Value = "Mobile No."
$ws.Range("B4").Value = "Email.ID"

# Give row 6, column B the built-in Hyperlink look (used later for the
# e-mail address column) before the shared-string table fills up so the
# "abc@gmail.com" / "9423944227" strings land in the same order as the
# source edit.
$ws.Range("B6").Style = "Hyperlink"

$ws.Range("B5").Value = "abc@gmail.com"
$ws.Range("A5").Value = "'9423944227"

# Row 6 stays empty but keeps formatting: quote-prefixed text style on A6,
# hyperlink style (already applied) on B6.
$ws.Range("A6").Value = "'1"
$ws.Range("A6").ClearContents()
$ws.Range("B6").ClearContents()

$ws.Rows("5:6").RowHeight = 15.75

# --- Cosmetic touch-ups ----------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 17.7

$ws.Range("B5").Select()
